# Resolve the ambiguous duplicate-name rows (RUT 21075353 shared by two
# students with similar names) by replacing them with a single, correctly
# selected student (RUT 21494146 - MAXIMILIANO JOAQUIN ALMONACID PEREZ),
# inserted as the new first data row, and removing the old duplicate rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2, shifting existing data rows down.
$ws.Rows.Item(2).Insert()

# Fill in the newly inserted row 2 with the selected student's data.
$ws.Cells.Item(2, 1).Value = "31-12-2024"
$ws.Cells.Item(2, 2).Value = "21494146"
$ws.Cells.Item(2, 3).Value = "5"
$ws.Cells.Item(2, 4).Value = "MAXIMILIANO JOAQUIN"
$ws.Cells.Item(2, 5).Value = "ALMONACID PÉREZ"
$ws.Cells.Item(2, 6).Value = "1"
$ws.Cells.Item(2, 7).Value = "FÍSICA MECANICA / 3"

# The former row 2 (MATIAS IGNACIO CEBALLOS VASQUEZ) is now row 3 - keep it.
# Remove the old ambiguous duplicate-name rows, which are now rows 4 and 5
# (FLAVIO ALEXANDER JARA LABRIN and ANA DORA LABRIN ESPINOZA, both RUT 21075353).
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()
